$d = $word.ActiveDocument

# The target paragraph is the last paragraph in the main body (the empty,
# centered paragraph that immediately precedes the final sectPr).
$last = $d.Paragraphs.Last

# 1) Drop the centered justification (w:jc="center") from the paragraph -
#    setting alignment back to left removes the <w:jc/> element entirely.
$last.Range.ParagraphFormat.Alignment = 0

# 2) Insert the hyperlink display text as two runs so that, once the
#    hyperlink style is applied, we end up with two <w:r> elements split
#    at "sumativa-" / "3" (matching the authored edit).
$insPt = $last.Range
$insPt.Collapse(0)
$startPos = $insPt.Start

$part1 = $d.Range($startPos, $startPos)
$part1.InsertAfter("https://github.com/flavioeichin/sumativa-")

$pos2 = $part1.End
$part2 = $d.Range($pos2, $pos2)
$part2.InsertAfter("3")

$fullRange = $d.Range($part1.Start, $part2.End)
$linkText = $fullRange.Text

# 3) Turn that text range into a real hyperlink pointing at the repo URL.
$h = $d.Hyperlinks.Add($fullRange, "https://github.com/flavioeichin/sumativa-3", $null, $null, $linkText)

# Make sure the run(s) use the document's existing "Hipervnculo" character
# style (the add above can otherwise stamp the style's display name
# instead of its style id).
$h.Range.Style = "Hipervnculo"

# Force the hyperlink text back into two separate runs (one per inserted
# chunk) by nudging formatting on just the last character ("3") and
# reverting it - this splits the run without altering final formatting.
$hStart = $h.Range.Start
$hEnd = $h.Range.End
$lastCharRange = $d.Range($hEnd - 2, $hEnd - 1)
$lastCharRange.Font.Bold = 1
$lastCharRange.Font.Bold = 0

# 4) Append a trailing literal space run right after the hyperlink.
$afterLink = $d.Range($hEnd - 1, $hEnd - 1)
$afterLink.InsertAfter(" ")
